# Insert a new data row at row 162 (pushing all existing rows 162..280 down to 163..281)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(162).Insert()

$ws.Cells.Item(162, 1).Value2 = 10
$ws.Cells.Item(162, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(162, 3).Value2 = "La Araucanía"
$ws.Cells.Item(162, 4).Value2 = 45072
$ws.Cells.Item(162, 5).Value2 = 9
$ws.Cells.Item(162, 6).Value2 = "Fruta"
$ws.Cells.Item(162, 7).Value2 = 100104
$ws.Cells.Item(162, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(162, 9).Value2 = 100104003
$ws.Cells.Item(162, 10).Value2 = "Membrillo"
$ws.Cells.Item(162, 11).Value2 = "Champion"
$ws.Cells.Item(162, 12).Value2 = "Primera"
$ws.Cells.Item(162, 13).Value2 = 145
$ws.Cells.Item(162, 14).Value2 = 14000
$ws.Cells.Item(162, 15).Value2 = 15000
$ws.Cells.Item(162, 16).Value2 = 14552
$ws.Cells.Item(162, 17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(162, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(162, 19).Value2 = 808
$ws.Cells.Item(162, 20).Value2 = 18
